$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header formatting (bold, centered, bordered) from an existing header cell
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Boolean outlier flags for data rows 2-25 (columns F=KNN, G=SVM, H=RF)
$knnFlags = @($False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False)
$svmFlags = @($False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $True,  $False, $False, $False, $False, $False, $False, $True,  $False, $False, $False)
$rfFlags  = @($False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False, $False)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $knnFlags[$i]
    $ws.Cells.Item($row, 7).Value = $svmFlags[$i]
    $ws.Cells.Item($row, 8).Value = $rfFlags[$i]
}
